# Update column G ("K") values on Sheet1 (rows 2-42) with regenerated
# strikeout counts (K) that replace the previous Strike# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 3
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 3
    12 = 1
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 1
    18 = 5
    19 = 2
    20 = 2
    21 = 2
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 2
    33 = 2
    34 = 2
    35 = 1
    36 = 2
    37 = 1
    38 = 3
    39 = 2
    40 = 0
    41 = 4
    42 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
